$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" values are stored as text in the sheet (scraped data),
# so force text format on the target cells before assigning the new values
# to avoid Excel auto-converting them to numbers (which would lose exact
# formatting such as trailing zeros / fixed-point notation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "259.37"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.59"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.133"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06109"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.561"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.509"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.332"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01328"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08131"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03525"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03185"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09211"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.771"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001641"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04662"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006428"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006128"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001070"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001503"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.732"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.265"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3315"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002718"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04601"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006999"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003706"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01014"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006091"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009917"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8040"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001128"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001904"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01242"
